# UndoRedoActivityDiagram.pptx update
#
# Mirrors the commit that reworks the undo/redo activity diagram to match
# the new UndoRedoCareTaker (ReadOnlyAddressBook list) design instead of
# the old UndoRedoStack / UndoableCommand design:
#   - removes the "[command is undoable]" branch (Diamond 11, the
#     undo/[else] text boxes, the "Clear redo stack" box and Diamond 25)
#   - slides the remaining flow shapes left/down to fill the gap
#   - renames/rewords the remaining decision + process boxes for the new
#     "address book state" based undo/redo flow

function Get-ShapeById {
    param($slide, $id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$EMU_PER_POINT = 12700

function Set-ShapeOffset {
    param($shape, $x, $y)
    $shape.Left = $x / $EMU_PER_POINT
    $shape.Top = $y / $EMU_PER_POINT
}

function Set-ShapeExtent {
    param($shape, $cx, $cy)
    $shape.Width = $cx / $EMU_PER_POINT
    $shape.Height = $cy / $EMU_PER_POINT
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Reposition the shapes that lead up to the removed decision branch
# ---------------------------------------------------------------------

$oval3 = Get-ShapeById $s 4
Set-ShapeOffset $oval3 1929588 3261938

$arrow5 = Get-ShapeById $s 6
Set-ShapeOffset $arrow5 2165257 3379773

$rect7 = Get-ShapeById $s 8
Set-ShapeOffset $rect7 2392862 3022393

$arrow8 = Get-ShapeById $s 9
Set-ShapeOffset $arrow8 3963217 3379774

# ---------------------------------------------------------------------
# 2. Remove the "[command is undoable]" branch and its connectors
# ---------------------------------------------------------------------

$idsToDelete = @(24, 55, 57, 61, 12, 19, 20, 22, 26)
foreach ($id in $idsToDelete) {
    $sh = Get-ShapeById $s $id
    if ($sh -ne $null) {
        $sh.Delete()
    }
}

# ---------------------------------------------------------------------
# 3. Reposition the remaining shapes to close the gap
# ---------------------------------------------------------------------

$diamond45 = Get-ShapeById $s 46
Set-ShapeOffset $diamond45 8215441 3248329

$textbox46 = Get-ShapeById $s 47
Set-ShapeOffset $textbox46 4801950 3488712

$textbox47 = Get-ShapeById $s 48
$textbox47.TextFrame.TextRange.Text = "[address book different from top of undo stack]"
Set-ShapeOffset $textbox47 2982124 2108748
Set-ShapeExtent $textbox47 2406969 646587

$rect50 = Get-ShapeById $s 51
$tr = $rect50.TextFrame.TextRange
$tr.Text = "Add "
$tr2 = $tr.InsertAfter("address book state ")
$tr3 = $tr2.InsertAfter("to undo ")
$tr4 = $tr3.InsertAfter("stack, clear redo stack")
Set-ShapeOffset $rect50 5389094 2377167
Set-ShapeExtent $rect50 2406970 888617

$diamond55 = Get-ShapeById $s 56
Set-ShapeOffset $diamond55 4356331 3140229

$arrow68 = Get-ShapeById $s 69
Set-ShapeOffset $arrow68 8696207 3488712

$group74 = Get-ShapeById $s 75
Set-ShapeOffset $group74 9115584 3378406

# ---------------------------------------------------------------------
# 4. Reposition the remaining elbow connectors
# ---------------------------------------------------------------------

$elbow65 = Get-ShapeById $s 66
Set-ShapeOffset $elbow65 4833528 2584663
Set-ShapeExtent $elbow65 318753 792380

$elbow71 = Get-ShapeById $s 72
$elbow71.Rotation = 270
$elbow71.VerticalFlip = 0
$elbow71.Adjustments.Item(1) = 3.11471
Set-ShapeOffset $elbow71 6472219 1745490
Set-ShapeExtent $elbow71 108100 3859110

$elbow73 = Get-ShapeById $s 74
Set-ShapeOffset $elbow73 7796064 2821476
Set-ShapeExtent $elbow73 659760 426853
